$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 183, shifting existing rows 183-215 down to 184-216.
$ws.Rows("183:183").Insert()

# Populate the newly inserted row 183 with the new record's data.
$ws.Range("A183").Value = 3
$ws.Range("B183").Value = "Femacal de La Calera"
$ws.Range("C183").Value = "Coquimbo"
$ws.Range("D183").Value = 44476
$ws.Range("D183").Style = "Normal"
$ws.Range("D183").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E183").Value = 5
$ws.Range("F183").Value = 100112040
$ws.Range("G183").Value = "Cilantro"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 160
$ws.Range("K183").Value = 2500
$ws.Range("L183").Value = 2500
$ws.Range("M183").Value = 2500
$ws.Range("N183").Value = "$/docena de atados (3 kilos)"
$ws.Range("O183").Value = "Provincia de Quillota"
$ws.Range("P183").Value = 833
$ws.Range("Q183").Value = 3
$ws.Range("R183").Value = "Hortaliza"
